# Updated cryptos list with refreshed prices / 1h volume percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds strings that often look numeric
# (e.g. "11.10", "1.010", "0.9502") and sometimes use "." as a thousands
# separator (e.g. "20.538.49"). Force the whole column to Text format
# before writing so Excel doesn't reinterpret the values as numbers and
# silently strip trailing zeros / collapse the thousand-dot grouping.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "20.538.49"
$ws.Range("E2").Value = "  +1.47%  "

$ws.Range("D3").Value = "1.474.61"
$ws.Range("E3").Value = "  +2.69%  "

$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").Value = "0.9502"
$ws.Range("E5").Value = "  +5.95%  "

$ws.Range("D6").Value = "278.19"
$ws.Range("E6").Value = "  +0.37%  "

$ws.Range("D7").Value = "0.3612"
$ws.Range("E7").Value = "  -1.52%  "

$ws.Range("D8").Value = "0.3051"
$ws.Range("E8").Value = "  -2.42%  "

$ws.Range("D9").Value = "39.44"
$ws.Range("E9").Value = "  +2.01%  "

$ws.Range("D10").Value = "1.055"
$ws.Range("E10").Value = "  +4.43%  "

$ws.Range("D11").Value = "0.06646"
$ws.Range("E11").Value = "  +1.89%  "

$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").Value = "5.506"
$ws.Range("E13").Value = "  +2.09%  "

$ws.Range("D14").Value = "18.04"
$ws.Range("E14").Value = "  +3.97%  "

$ws.Range("D15").Value = "6.194"
$ws.Range("E15").Value = "  +1.89%  "

$ws.Range("D16").Value = "0.9509"
$ws.Range("E16").Value = "  +5.45%  "

$ws.Range("D17").Value = "0.00001027"
$ws.Range("E17").Value = "  +0.83%  "

$ws.Range("D18").Value = "1.472.38"
$ws.Range("E18").Value = "  +2.36%  "

$ws.Range("D19").Value = "0.05934"
$ws.Range("E19").Value = "  +5.73%  "

$ws.Range("D20").Value = "69.29"
$ws.Range("E20").Value = "  +3.72%  "

$ws.Range("D21").Value = "5.488"
$ws.Range("E21").Value = "  +0.49%  "

$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("D23").Value = "11.10"
$ws.Range("E23").Value = "  +1.49%  "

$ws.Range("D24").Value = "2.263"
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").Value = "20.593.52"
$ws.Range("E25").Value = "  +1.44%  "

$ws.Range("D26").Value = "142.98"
$ws.Range("E26").Value = "  +5.47%  "

$ws.Range("D27").Value = "2.114"
$ws.Range("E27").Value = "  -3.31%  "

$ws.Range("D28").Value = "17.18"
$ws.Range("E28").Value = "  +1.61%  "

$ws.Range("D29").Value = "1.631.58"
$ws.Range("E29").Value = "  +2.15%  "

$ws.Range("D30").Value = "113.58"
$ws.Range("E30").Value = "  +2.66%  "

$ws.Range("D31").Value = "3.949"
$ws.Range("E31").Value = "  +10.08%  "

$ws.Range("D32").Value = "5.004"
$ws.Range("E32").Value = "  +2.62%  "

$ws.Range("D33").Value = "0.8076"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").Value = "0.07963"
$ws.Range("E34").Value = "  +4.43%  "

$ws.Range("D35").Value = "1.510"
$ws.Range("E35").Value = "  +4.72%  "

$ws.Range("D36").Value = "1.223"
$ws.Range("E36").Value = "  +8.64%  "

$ws.Range("D37").Value = "0.05845"
$ws.Range("E37").Value = "  -1.32%  "

$ws.Range("D38").Value = "4.724"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("D39").Value = "0.02042"
$ws.Range("E39").Value = "  +1.97%  "

$ws.Range("D40").Value = "10.33"
$ws.Range("E40").Value = "  +0.63%  "

$ws.Range("E41").Value = "  +4.03%  "

$ws.Range("D42").Value = "0.1872"
$ws.Range("E42").Value = "  +2.58%  "

$ws.Range("D43").Value = "7.419"
$ws.Range("E43").Value = "  +9.98%  "

$ws.Range("D44").Value = "0.5282"
$ws.Range("E44").Value = "  +0.78%  "

# Rows 45/46 swap places: EnergySwap now ranks above PancakeSwap, with
# refreshed price / volume figures for both.
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "12.31"
$ws.Range("E45").Value = "  +2.82%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.528"
$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("D47").Value = "117.89"
$ws.Range("E47").Value = "  -1.45%  "

$ws.Range("D48").Value = "0.5185"
$ws.Range("E48").Value = "  +1.03%  "

$ws.Range("D49").Value = "1.812"
$ws.Range("E49").Value = "  +3.14%  "

$ws.Range("D50").Value = "0.06477"
$ws.Range("E50").Value = "  +3.00%  "

$ws.Range("D51").Value = "0.9795"
$ws.Range("E51").Value = "  -2.16%  "
